$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.835941000000001
$ws.Range("H2").Value = 17.507823
$ws.Range("I2").Value = 0.03643643319117328
$ws.Range("J2").Value = 0.03643643319117327
$ws.Range("M2").Value = 0.24449
$ws.Range("N2").Value = 0.73347
$ws.Range("O2").Value = 0.009675524511058336
$ws.Range("P2").Value = 0.009675524511058334
$ws.Range("Q2").Value = 1.42682921509
$ws.Range("R2").Value = 12.84146293581
$ws.Range("S2").Value = 0.0003525416024367365
$ws.Range("T2").Value = 0.0003525416024367364

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.835941000000001
$ws.Range("H3").Value = 17.507823
$ws.Range("I3").Value = 0.03643643319117328
$ws.Range("J3").Value = 0.03643643319117327
$ws.Range("O3").Value = 0.1748614367985708
$ws.Range("P3").Value = 0.1748614367985708
$ws.Range("Q3").Value = 25.78644768370534
$ws.Range("R3").Value = 232.078029153348
$ws.Range("S3").Value = 0.006371327059623694
$ws.Range("T3").Value = 0.006371327059623691

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.835941000000001
$ws.Range("H4").Value = 17.507823
$ws.Range("I4").Value = 0.03643643319117328
$ws.Range("J4").Value = 0.03643643319117327
$ws.Range("M4").Value = 20.60586566666667
$ws.Range("N4").Value = 61.81759700000001
$ws.Range("O4").Value = 0.815463038690371
$ws.Range("P4").Value = 0.8154630386903708
$ws.Range("Q4").Value = 120.2546162845924
$ws.Range("R4").Value = 1082.291546561331
$ws.Range("S4").Value = 0.02971256452911285
$ws.Range("T4").Value = 0.02971256452911284

# Row 5
$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.1093102818770573
$ws.Range("J5").Value = 0.1093102818770573
$ws.Range("M5").Value = 0.24449
$ws.Range("N5").Value = 0.73347
$ws.Range("O5").Value = 0.009675524511058336
$ws.Range("P5").Value = 0.009675524511058334
$ws.Range("Q5").Value = 4.280526111696667
$ws.Range("R5").Value = 38.52473500527
$ws.Range("S5").Value = 0.001057634311612164
$ws.Range("T5").Value = 0.001057634311612163

# Row 6
$ws.Range("G6").Value = 17.50798033333334
$ws.Range("H6").Value = 52.52394100000001
$ws.Range("I6").Value = 0.1093102818770573
$ws.Range("J6").Value = 0.1093102818770573
$ws.Range("O6").Value = 0.1748614367985708
$ws.Range("P6").Value = 0.1748614367985708
$ws.Range("Q6").Value = 77.36003823767957
$ws.Range("R6").Value = 696.240344139116
$ws.Range("S6").Value = 0.01911415294587901
$ws.Range("T6").Value = 0.01911415294587901

# Row 7
$ws.Range("G7").Value = 17.50798033333334
$ws.Range("H7").Value = 52.52394100000001
$ws.Range("I7").Value = 0.1093102818770573
$ws.Range("J7").Value = 0.1093102818770573
$ws.Range("M7").Value = 20.60586566666667
$ws.Range("N7").Value = 61.81759700000001
$ws.Range("O7").Value = 0.815463038690371
$ws.Range("P7").Value = 0.8154630386903708
$ws.Range("Q7").Value = 360.7670908433087
$ws.Range("R7").Value = 3246.903817589778
$ws.Range("S7").Value = 0.08913849461956613
$ws.Range("T7").Value = 0.08913849461956611

# Row 8
$ws.Range("G8").Value = 136.8238143333333
$ws.Range("H8").Value = 410.471443
$ws.Range("I8").Value = 0.8542532849317694
$ws.Range("J8").Value = 0.8542532849317694
$ws.Range("M8").Value = 0.24449
$ws.Range("N8").Value = 0.73347
$ws.Range("O8").Value = 0.009675524511058336
$ws.Range("P8").Value = 0.009675524511058334
$ws.Range("Q8").Value = 33.45205436635666
$ws.Range("R8").Value = 301.06848929721
$ws.Range("S8").Value = 0.008265348597009435
$ws.Range("T8").Value = 0.008265348597009434

# Row 9
$ws.Range("G9").Value = 136.8238143333333
$ws.Range("H9").Value = 410.471443
$ws.Range("I9").Value = 0.8542532849317694
$ws.Range("J9").Value = 0.8542532849317694
$ws.Range("O9").Value = 0.1748614367985708
$ws.Range("P9").Value = 0.1748614367985708
$ws.Range("Q9").Value = 604.5640506289409
$ws.Range("R9").Value = 5441.076455660468
$ws.Range("S9").Value = 0.1493759567930681
$ws.Range("T9").Value = 0.1493759567930681

# Row 10
$ws.Range("G10").Value = 136.8238143333333
$ws.Range("H10").Value = 410.471443
$ws.Range("I10").Value = 0.8542532849317694
$ws.Range("J10").Value = 0.8542532849317694
$ws.Range("M10").Value = 20.60586566666667
$ws.Range("N10").Value = 61.81759700000001
$ws.Range("O10").Value = 0.815463038690371
$ws.Range("P10").Value = 0.8154630386903708
$ws.Range("Q10").Value = 2819.373138153608
$ws.Range("R10").Value = 25374.35824338247
$ws.Range("S10").Value = 0.696611979541692
$ws.Range("T10").Value = 0.6966119795416919
